# Remove the slide-number placeholder shape from the first slide.
# (Turning off the slide's slide-number display removes the
# "Slide Number Placeholder 1" <p:sp> that PowerPoint had instantiated
# from the layout's sldNum placeholder.)
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.HeadersFooters.SlideNumber.Visible = $false
